# Applies the Mar 30 2024 08:24:23 UTC cryptos-list refresh (prices, 1h volume %,
# and the Dai/PEPE + Stellar/ApeXProtocol row-content swaps) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = '69.837.17'
$ws.Range("E2").Value = '  +0.02%  '

# Row 3 (Ethereum)
$ws.Range("D3").Value = '3.502.18'
$ws.Range("E3").Value = '  -0.72%  '

# Row 4 (TetherUSD)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.53%  '

# Row 6 (Solana)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '195.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.81%  '

# Row 7 (XRP)
$ws.Range("E7").Value = '  +1.91%  '

# Row 8 (USDC)
$ws.Range("E8").Value = '  -0.10%  '

# Row 9 (Dogecoin)
$ws.Range("E9").Value = '  -2.02%  '

# Row 10 (Cardano)
$ws.Range("E10").Value = '  +1.79%  '

# Row 11 (Avalanche)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '

# Row 12 (ShibaInu)
$ws.Range("E12").Value = '  -2.04%  '

# Row 13 (Polkadot)
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.36%  '

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = '4.053.68'
$ws.Range("E14").Value = '  -0.84%  '

# Row 15 (BitcoinCash)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '603.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.03%  '

# Row 16 (WrappedBTC)
$ws.Range("D16").Value = '69.941.25'
$ws.Range("E16").Value = '  +0.15%  '

# Row 17 (Chainlink)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.23%  '

# Row 18 (Uniswap)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.02%  '

# Row 19 (WrappedEther)
$ws.Range("D19").Value = '3.500.06'
$ws.Range("E19").Value = '  -1.53%  '

# Row 20 (TRON)
$ws.Range("E20").Value = '  +0.64%  '

# Row 21 (Polygon)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.992'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.46%  '

# Row 22 (InternetComputer(DFINITY))
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.02%  '

# Row 23 (Litecoin)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '104.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +11.08%  '

# Row 24 (Toncoin)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.11%  '

# Row 25 (PancakeSwap)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.45%  '

# Row 26 (ImmutableX)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.95%  '

# Row 27 (RenderToken)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.22%  '

# Row 28 (Filecoin)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.23%  '

# Row 29 (EthereumClassic)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.78%  '

# Row 30 (dogwifhat)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +26.12%  '

# Row 31 (NEARProtocol)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.94%  '

# Row 32 (Cosmos)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.31%  '

# Row 33 (Hedera)
$ws.Range("E33").Value = '  +1.70%  '

# Row 34 (OKB)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.80%  '

# Row 35 (Maker)
$ws.Range("D35").Value = '3.737.32'
$ws.Range("E35").Value = '  +6.30%  '

# Row 36 (Dai)
$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").Value = '0.0₃0810'
$ws.Range("E36").Value = '  +4.37%  '

# Row 37 (PEPE)
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.23%  '

# Row 38 (Fetch.AI)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.24%  '

# Row 39 (TheGraph)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.391'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.66%  '

# Row 40 (InjectiveProtocol)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.70'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.71%  '

# Row 41 (Bittensor)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '503.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.32%  '

# Row 42 (Stacks)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.04%  '

# Row 43 (Kaspa)
$ws.Range("E43").Value = '  +0.59%  '

# Row 44 (VeChain)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0458'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.32%  '

# Row 45 (ApeXProtocol)
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.140'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.23%  '

# Row 46 (Stellar)
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.41%  '

# Row 47 (ThetaToken)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.04%  '

# Row 48 (FirstDigitalUSD)
$ws.Range("E48").Value = '  +0.40%  '

# Row 49 (THORChain)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.28%  '

# Row 50 (Monero)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.74'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.89%  '

# Row 51 (FLOKI)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000242'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.62%  '
